$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6 (diff @ old line 917)
$ws.Range("H6").Value = 388.0909
$ws.Range("J6").Value = 644.75
$ws.Range("L6").Value = 1934.25
$ws.Range("N6").Value = -2158.25

# Row 7 (diff @ old line 969)
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").ClearContents()
$ws.Range("N7").Value = 0

# Row 14 (diff @ old line 1318)
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").ClearContents()
$ws.Range("N14").Value = 0

# Row 15 (diff @ old line 1367)
$ws.Range("H15").Value = 313.82224
$ws.Range("I15").Value = 313.82224
$ws.Range("K15").Value = 941.4667200000001
$ws.Range("M15").Value = -772.4667200000001

# Row 31 (diff @ old line 2163)
$ws.Range("H31").Value = 751.5
$ws.Range("I31").Value = 751.5
$ws.Range("K31").Value = 2254.5
$ws.Range("M31").Value = -2024.5

# Row 32 (diff @ old line 2212)
$ws.Range("H32").Value = 836.5
$ws.Range("J32").Value = 919.0714
$ws.Range("L32").Value = 919.0714
$ws.Range("N32").Value = -1571.0714

# Row 38 (diff @ old line 2518)
$ws.Range("H38").Value = 4927
$ws.Range("I38").Value = 228.6
$ws.Range("J38").Value = 10800
$ws.Range("K38").Value = 685.8
$ws.Range("L38").Value = 32400
$ws.Range("M38").Value = -313.8
$ws.Range("N38").Value = -33144

# Row 39 (diff @ old line 2570)
$ws.Range("H39").Value = 225.33333
$ws.Range("I39").Value = 115
$ws.Range("J39").Value = 666.6667
$ws.Range("K39").Value = 345
$ws.Range("L39").Value = 2000.0001
$ws.Range("M39").Value = -49
$ws.Range("N39").Value = -2592.0001

# Row 103 (diff @ old line 5793)
$ws.Range("H103").Value = 4845.619
$ws.Range("J103").Value = 11425
$ws.Range("L103").Value = 34275
$ws.Range("N103").Value = -35447

# Row 106 (diff @ old line 5940)
$ws.Range("H106").Value = 2563.4285
$ws.Range("I106").Value = 1853.6364
$ws.Range("J106").Value = 5166
$ws.Range("K106").Value = 1853.6364
$ws.Range("L106").Value = 5166
$ws.Range("M106").Value = -1222.6364
$ws.Range("N106").Value = -6428

# Row 113 (diff @ old line 6289)
$ws.Range("H113").Value = 7935.36
$ws.Range("I113").Value = 5155.25
$ws.Range("J113").Value = 12877.777
$ws.Range("K113").Value = 5155.25
$ws.Range("L113").Value = 12877.777
$ws.Range("M113").Value = -1901.25
$ws.Range("N113").Value = -19385.777

# Row 123 (diff @ old line 6782)
$ws.Range("H123").Value = 41797.777
$ws.Range("J123").Value = 41797.777
$ws.Range("L123").Value = 41797.777
$ws.Range("N123").Value = -51597.777

# Row 138 (diff @ old line 7535)
$ws.Range("H138").Value = 2348.1025
$ws.Range("I138").Value = 1197
$ws.Range("J138").Value = 3237.5908
$ws.Range("K138").Value = 3591
$ws.Range("L138").Value = 9712.7724
$ws.Range("M138").Value = 1549
$ws.Range("N138").Value = -19992.7724


$ws = $wb.Worksheets.Item("ARM")
# Row 129 (diff @ old line 14027)
$ws.Range("H129").Value = 49780
$ws.Range("J129").Value = 49780
$ws.Range("L129").Value = 49780
$ws.Range("N129").Value = -59780

# Row 132 (diff @ old line 14174)
$ws.Range("H132").Value = 2690.848
$ws.Range("I132").Value = 1802.4073
$ws.Range("J132").Value = 3953.3684
$ws.Range("K132").Value = 5407.2219
$ws.Range("L132").Value = 11860.1052
$ws.Range("M132").Value = -2877.2219
$ws.Range("N132").Value = -16920.1052


$ws = $wb.Worksheets.Item("BSM")
# Row 20 (diff @ old line 15640)
$ws.Range("H20").Value = 10033.053
$ws.Range("I20").Value = 1370.8889
$ws.Range("J20").Value = 17829
$ws.Range("K20").Value = 1370.8889
$ws.Range("L20").Value = 17829
$ws.Range("M20").Value = -1123.8889
$ws.Range("N20").Value = -18323


$ws = $wb.Worksheets.Item("CRP")
# Row 31 (diff @ old line 23133)
$ws.Range("H31").Value = 189623.62
$ws.Range("I31").Value = 451213.7
$ws.Range("J31").Value = 2773.5715
$ws.Range("K31").Value = 451213.7
$ws.Range("L31").Value = 2773.5715
$ws.Range("M31").Value = -450918.7
$ws.Range("N31").Value = -3363.5715

# Row 34 (diff @ old line 23289)
$ws.Range("H34").Value = 189623.62
$ws.Range("I34").Value = 451213.7
$ws.Range("J34").Value = 2773.5715
$ws.Range("K34").Value = 451213.7
$ws.Range("L34").Value = 2773.5715
$ws.Range("M34").Value = -451011.7
$ws.Range("N34").Value = -3177.5715

# Row 99 (diff @ old line 26507)
$ws.Range("H99").Value = 14290380
$ws.Range("I99").Value = 22224502
$ws.Range("J99").Value = 8960
$ws.Range("K99").Value = 22224502
$ws.Range("L99").Value = 8960
$ws.Range("M99").Value = -22223004
$ws.Range("N99").Value = -11956

# Row 126 (diff @ old line 27824)
$ws.Range("H126").Value = 14290380
$ws.Range("I126").Value = 22224502
$ws.Range("J126").Value = 8960
$ws.Range("K126").Value = 66673506
$ws.Range("L126").Value = 26880
$ws.Range("M126").Value = -66671036
$ws.Range("N126").Value = -31820

# Row 134 (diff @ old line 28222)
$ws.Range("H134").Value = 5816.6665
$ws.Range("I134").Value = 1225
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 3675
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -1140
$ws.Range("N134").Value = -50070


$ws = $wb.Worksheets.Item("CUL")
# Row 33 (diff @ old line 30281)
$ws.Range("H33").Value = 131.23077
$ws.Range("I33").Value = 107.36364
$ws.Range("J33").Value = 148.73334
$ws.Range("K33").Value = 644.18184
$ws.Range("L33").Value = 892.40004
$ws.Range("M33").Value = -361.18184
$ws.Range("N33").Value = -1458.40004

# Row 37 (diff @ old line 30477)
$ws.Range("H37").Value = 73444.44500000001
$ws.Range("J37").Value = 73444.44500000001
$ws.Range("L37").Value = 220333.335
$ws.Range("N37").Value = -220557.335

# Row 64 (diff @ old line 31848)
$ws.Range("H64").Value = 7000
$ws.Range("I64").Value = 2600
$ws.Range("J64").Value = 9200
$ws.Range("K64").Value = 7800
$ws.Range("L64").Value = 27600
$ws.Range("M64").Value = -7530
$ws.Range("N64").Value = -28140

# Row 67 (diff @ old line 32001)
$ws.Range("H67").Value = 7000
$ws.Range("I67").Value = 2600
$ws.Range("J67").Value = 9200
$ws.Range("K67").Value = 7800
$ws.Range("L67").Value = 27600
$ws.Range("M67").Value = -6864
$ws.Range("N67").Value = -29472

# Row 68 (diff @ old line 32053)
$ws.Range("H68").Value = 1026.4531
$ws.Range("I68").Value = 733.85187
$ws.Range("J68").Value = 1239.973
$ws.Range("K68").Value = 2201.55561
$ws.Range("L68").Value = 3719.919
$ws.Range("M68").Value = -1390.55561
$ws.Range("N68").Value = -5341.919

# Row 71 (diff @ old line 32209)
$ws.Range("H71").Value = 1026.4531
$ws.Range("I71").Value = 733.85187
$ws.Range("J71").Value = 1239.973
$ws.Range("K71").Value = 6604.66683
$ws.Range("L71").Value = 11159.757
$ws.Range("M71").Value = -2548.66683
$ws.Range("N71").Value = -19271.757

# Row 113 (diff @ old line 34333)
$ws.Range("H113").Value = 539.05554
$ws.Range("I113").Value = 463.84
$ws.Range("J113").Value = 603.89655
$ws.Range("K113").Value = 1391.52
$ws.Range("L113").Value = 1811.68965
$ws.Range("M113").Value = 778.48
$ws.Range("N113").Value = -6151.68965

# Row 131 (diff @ old line 35248)
$ws.Range("H131").Value = 780.5454999999999
$ws.Range("J131").Value = 801.44086
$ws.Range("L131").Value = 2404.32258
$ws.Range("N131").Value = -12484.32258


$ws = $wb.Worksheets.Item("GSM")
# Row 70 (diff @ old line 39252)
$ws.Range("H70").Value = 5872.5894
$ws.Range("I70").Value = 5583.136
$ws.Range("J70").Value = 6933.9165
$ws.Range("K70").Value = 5583.136
$ws.Range("L70").Value = 6933.9165
$ws.Range("M70").Value = -5313.136
$ws.Range("N70").Value = -7473.9165

# Row 73 (diff @ old line 39399)
$ws.Range("H73").Value = 5872.5894
$ws.Range("I73").Value = 5583.136
$ws.Range("J73").Value = 6933.9165
$ws.Range("K73").Value = 5583.136
$ws.Range("L73").Value = 6933.9165
$ws.Range("M73").Value = -4647.136
$ws.Range("N73").Value = -8805.916499999999

# Row 123 (diff @ old line 41843)
$ws.Range("H123").Value = 10485.25
$ws.Range("J123").Value = 10485.25
$ws.Range("L123").Value = 10485.25
$ws.Range("N123").Value = -15385.25

